$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Holly added "S.GISH" to the harvester column (column B) in bioSamples;
# reflect the same fix here in rnaSamples: update the harvester values
# (B2:B24) from "Retrofitted_1344" to "S.GISH".
$lastRow = $ws.Cells.Item($ws.Rows.Count, 2).End(-4162).Row
if ($lastRow -lt 2) { $lastRow = 24 }

for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 2).Value = "S.GISH"
}

# Mirror the resulting selection on column B, as left by the edit.
$ws.Range("B:B").Select() | Out-Null

$wb.Save()
